$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for 2 (44687/160) and 4 (44691/100) swap their
# Date (D) and Volumen (J) values.
$ws.Range("D2").Value = 44691
$ws.Range("J2").Value = 100

$ws.Range("D4").Value = 44687
$ws.Range("J4").Value = 160

# The data rows for 3 and 5 swap their Date (D), Volumen (J),
# Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M),
# Unidad de comercializacion (N), Origen (O), Precio $/Kg (P) and
# Kg o Unidades (Q) values.
$ws.Range("D3").Value = 44692
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 3500
$ws.Range("M3").Value = 3250
$ws.Range("N3").Value = "$/docena de matas"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 542
$ws.Range("Q3").Value = 6

$ws.Range("D5").Value = 44221
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 1300
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1420
$ws.Range("N5").Value = "$/atado"
$ws.Range("O5").Value = "Provincia de Diguillín"
$ws.Range("P5").Value = 1420
$ws.Range("Q5").Value = 1
